# Opdateret tidsplan for den 14-03-2017 jeppe
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tidsregistrering")

# ---------------------------------------------------------------------------
# New day entry: 14-03-2017 (serial 42808) -> rows 29-33
# Formatting is copied from the most similar existing cell in each column so
# the new cells land on the same style indices Excel itself would reuse.
# ---------------------------------------------------------------------------

# Row 29 - Test Analyst, 08:05 - 11:00
$ws.Range("A24").Copy()
$ws.Range("A29").PasteSpecial(-4122)
$ws.Range("A29").Value = 42808

$ws.Range("E27").Copy()
$ws.Range("E29").PasteSpecial(-4122)
$ws.Range("E29").Value = "Test Analyst"

$ws.Range("F24").Copy()
$ws.Range("F29").PasteSpecial(-4122)
$ws.Range("F29").Value = "Lavet brugertest på alpha brugergrænseflade"

$ws.Range("G24").Copy()
$ws.Range("G29").PasteSpecial(-4122)
$ws.Range("G29").Value = 0.33680555555555558

$ws.Range("H24").Copy()
$ws.Range("H29").PasteSpecial(-4122)
$ws.Range("H29").Value = 0.45833333333333331

# Row 30 - Requirements Specifier, 12:10 - 13:30
$ws.Range("E21").Copy()
$ws.Range("E30").PasteSpecial(-4122)
$ws.Range("E30").Value = "Requirements Specifier"

$ws.Range("F25").Copy()
$ws.Range("F30").PasteSpecial(-4122)
$ws.Range("F30").Value = "Lavet OC 5 beregnSigmaN"

$ws.Range("G25").Copy()
$ws.Range("G30").PasteSpecial(-4122)
$ws.Range("G30").Value = 0.50694444444444442

$ws.Range("H25").Copy()
$ws.Range("H30").PasteSpecial(-4122)
$ws.Range("H30").Value = 0.5625

# Row 32 text is entered before row 31's so the new shared-string table keeps
# the same ordering (OC5 Test Suite before OC7 Test Suite) as the source file.
$ws.Range("F27").Copy()
$ws.Range("F32").PasteSpecial(-4122)
$ws.Range("F32").Value = "Lavet Test Suite for OC5 beregnSigmaN"

# Row 31 - 13:35 - 14:04
$ws.Range("F26").Copy()
$ws.Range("F31").PasteSpecial(-4122)
$ws.Range("F31").Value = "Lavet Test Suite for OC7 beregnSigmaTau"

$ws.Range("G26").Copy()
$ws.Range("G31").PasteSpecial(-4122)
$ws.Range("G31").Value = 0.56597222222222221

$ws.Range("H26").Copy()
$ws.Range("H31").PasteSpecial(-4122)
$ws.Range("H31").Value = 0.58611111111111114

# Row 32 (remaining cells) - 14:45 - 14:50
$ws.Range("G27").Copy()
$ws.Range("G32").PasteSpecial(-4122)
$ws.Range("G32").Value = 0.61458333333333337

$ws.Range("H27").Copy()
$ws.Range("H32").PasteSpecial(-4122)
$ws.Range("H32").Value = 0.61805555555555558

# Row 33 - daily total
$ws.Range("I28").Copy()
$ws.Range("I33").PasteSpecial(-4122)
$ws.Range("I33").Value = 5

$excel.CutCopyMode = 0

# View state: scroll down a bit and leave selection where data entry ended
$ws.Application.ActiveWindow.ScrollRow = 7
$ws.Range("H34").Select()
